$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove row 6 (dataset now has one fewer sample row; dimension becomes A1:AH5)
$ws.Rows.Item(6).Delete()

# Column width tweaks (auto-fit style adjustments caused by the new data values)
$ws.Columns.Item(11).ColumnWidth = 6.17
$ws.Columns.Item(24).ColumnWidth = 7.17
$ws.Columns.Item(27).ColumnWidth = 7.17
$ws.Columns.Item(28).ColumnWidth = 7.17
$ws.Columns.Item(29).ColumnWidth = 7.17

# Row 2
$ws.Range("A2").Value = 45181.50694444445
$ws.Range("B2").Value = 19.217
$ws.Range("C2").Value = 12.901
$ws.Range("D2").Value = 4.042
$ws.Range("E2").Value = 40.812
$ws.Range("F2").Value = 32.818
$ws.Range("G2").Value = 15.123
$ws.Range("H2").Value = 47.986
$ws.Range("I2").Value = 23.269
$ws.Range("J2").Value = 9.710000000000001
$ws.Range("K2").Value = 14.67
$ws.Range("L2").Value = 16.076
$ws.Range("M2").Value = 16.742
$ws.Range("N2").Value = 4.827
$ws.Range("O2").Value = 15.038
$ws.Range("P2").Value = 20.994
$ws.Range("Q2").Value = 12.85
$ws.Range("R2").Value = 3.46
$ws.Range("S2").Value = 2.249
$ws.Range("T2").Value = 221.547
$ws.Range("U2").Value = 41.81
$ws.Range("V2").Value = 13.881
$ws.Range("W2").Value = 27.553
$ws.Range("X2").Value = 14.055
$ws.Range("Y2").Value = 3.03
$ws.Range("Z2").Value = 24.312
$ws.Range("AA2").Value = 12.261
$ws.Range("AB2").Value = 11.125
$ws.Range("AC2").Value = 13.047
$ws.Range("AD2").Value = 16.565
$ws.Range("AE2").Value = 3.456
$ws.Range("AF2").Value = 42.557
$ws.Range("AG2").Value = 7.647
$ws.Range("AH2").Value = 17.354

# Row 3
$ws.Range("A3").Value = 45181.51388888889
$ws.Range("B3").Value = 12.971
$ws.Range("C3").Value = 9.032
$ws.Range("D3").Value = 1.695
$ws.Range("E3").Value = 27.971
$ws.Range("F3").Value = 22.532
$ws.Range("G3").Value = 10.208
$ws.Range("H3").Value = 40.423
$ws.Range("I3").Value = 15.707
$ws.Range("J3").Value = 6.724
$ws.Range("K3").Value = 9.909000000000001
$ws.Range("L3").Value = 11.183
$ws.Range("M3").Value = 11.647
$ws.Range("N3").Value = 3.262
$ws.Range("O3").Value = 10.151
$ws.Range("P3").Value = 14.263
$ws.Range("Q3").Value = 8.82
$ws.Range("R3").Value = 1.526
$ws.Range("S3").Value = 0.963
$ws.Range("T3").Value = 147.193
$ws.Range("U3").Value = 28.454
$ws.Range("V3").Value = 9.369999999999999
$ws.Range("W3").Value = 18.788
$ws.Range("X3").Value = 9.728
$ws.Range("Y3").Value = 1.908
$ws.Range("Z3").Value = 19.494
$ws.Range("AA3").Value = 8.276
$ws.Range("AB3").Value = 7.538
$ws.Range("AC3").Value = 8.83
$ws.Range("AD3").Value = 11.56
$ws.Range("AE3").Value = 1.265
$ws.Range("AF3").Value = 36.84
$ws.Range("AG3").Value = 5.152
$ws.Range("AH3").Value = 11.715

# Row 4
$ws.Range("A4").Value = 45181.52083333334
$ws.Range("B4").Value = 6.726
$ws.Range("C4").Value = 4.589
$ws.Range("D4").Value = 0.997
$ws.Range("E4").Value = 14.515
$ws.Range("F4").Value = 11.55
$ws.Range("G4").Value = 5.293
$ws.Range("H4").Value = 23.321
$ws.Range("I4").Value = 8.144
$ws.Range("J4").Value = 3.457
$ws.Range("K4").Value = 4.989
$ws.Range("L4").Value = 5.824
$ws.Range("M4").Value = 6.03
$ws.Range("N4").Value = 1.694
$ws.Range("O4").Value = 5.263
$ws.Range("P4").Value = 7.372
$ws.Range("Q4").Value = 4.686
$ws.Range("R4").Value = 0.969
$ws.Range("S4").Value = 0.532
$ws.Range("T4").Value = 72.80500000000001
$ws.Range("U4").Value = 14.85
$ws.Range("V4").Value = 4.858
$ws.Range("W4").Value = 9.718999999999999
$ws.Range("X4").Value = 5.014
$ws.Range("Y4").Value = 1.09
$ws.Range("Z4").Value = 10.892
$ws.Range("AA4").Value = 4.291
$ws.Range("AB4").Value = 3.965
$ws.Range("AC4").Value = 4.633
$ws.Range("AD4").Value = 5.999
$ws.Range("AE4").Value = 0.773
$ws.Range("AF4").Value = 21.32
$ws.Range("AG4").Value = 2.622
$ws.Range("AH4").Value = 6.075

# Row 5
$ws.Range("A5").Value = 45181.52777777778
$ws.Range("B5").Value = 3.84
$ws.Range("C5").Value = 2.54
$ws.Range("D5").Value = 0.68
$ws.Range("E5").Value = 8.289999999999999
$ws.Range("F5").Value = 6.49
$ws.Range("G5").Value = 3.03
$ws.Range("H5").Value = 14.01
$ws.Range("I5").Value = 4.65
$ws.Range("J5").Value = 1.95
$ws.Range("K5").Value = 2.74
$ws.Range("L5").Value = 3.34
$ws.Range("M5").Value = 3.43
$ws.Range("N5").Value = 0.97
$ws.Range("O5").Value = 3.01
$ws.Range("P5").Value = 4.19
$ws.Range("Q5").Value = 2.76
$ws.Range("R5").Value = 0.71
$ws.Range("S5").Value = 0.34
$ws.Range("T5").Value = 38.47
$ws.Range("U5").Value = 8.529999999999999
$ws.Range("V5").Value = 2.78
$ws.Range("W5").Value = 5.52
$ws.Range("X5").Value = 2.84
$ws.Range("Y5").Value = 0.7
$ws.Range("Z5").Value = 6.45
$ws.Range("AA5").Value = 2.45
$ws.Range("AB5").Value = 2.31
$ws.Range("AC5").Value = 2.69
$ws.Range("AD5").Value = 3.42
$ws.Range("AE5").Value = 0.5600000000000001
$ws.Range("AF5").Value = 12.82
$ws.Range("AG5").Value = 1.46
$ws.Range("AH5").Value = 3.47
